# Table4.xlsx — minor edits: add "Sacramento" label to row 3, adjust a few
# row heights, and move the active-cell selection to A4.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A3 was blank; it now carries the "Sacramento" row label (adds a new shared
# string entry automatically).
$ws.Range("A3").Value = "Sacramento"

# A handful of rows were resized.
$ws.Range("A1").EntireRow.RowHeight = 16
$ws.Range("A3").EntireRow.RowHeight = 16
$ws.Range("A4").EntireRow.RowHeight = 17
$ws.Range("A5").EntireRow.RowHeight = 17
$ws.Range("A6").EntireRow.RowHeight = 16
$ws.Range("A7").EntireRow.RowHeight = 17
$ws.Range("A10").EntireRow.RowHeight = 17
$ws.Range("A12").EntireRow.RowHeight = 17

# The active selection moved from B12 to A4.
$ws.Range("A4").Select() | Out-Null
